$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "46.380.10"
Set-TextValue $ws "E2" "  -0.52%  "
Set-TextValue $ws "D3" "2.617.23"
Set-TextValue $ws "E3" "  +1.25%  "
Set-TextValue $ws "E4" "  +0.03%  "
Set-TextValue $ws "D5" "308.67"
Set-TextValue $ws "E5" "  +1.15%  "
Set-TextValue $ws "D6" "100.00"
Set-TextValue $ws "E6" "  -2.26%  "
Set-TextValue $ws "E7" "  +0.08%  "
Set-TextValue $ws "E8" "  +0.03%  "
Set-TextValue $ws "D9" "0.582"
Set-TextValue $ws "E9" "  +1.87%  "
Set-TextValue $ws "D10" "39.34"
Set-TextValue $ws "E10" "  +1.75%  "
Set-TextValue $ws "D11" "0.0846"
Set-TextValue $ws "E11" "  +1.59%  "
Set-TextValue $ws "D12" "54.34"
Set-TextValue $ws "E12" "  -1.11%  "
Set-TextValue $ws "D13" "8.16"
Set-TextValue $ws "E13" "  +1.86%  "
Set-TextValue $ws "D14" "3.014.67"
Set-TextValue $ws "E14" "  +1.40%  "
Set-TextValue $ws "E15" "  +0.77%  "
Set-TextValue $ws "D16" "2.615.31"
Set-TextValue $ws "E16" "  +0.74%  "
Set-TextValue $ws "D17" "0.927"
Set-TextValue $ws "E17" "  +3.24%  "
Set-TextValue $ws "E18" "  -0.15%  "
Set-TextValue $ws "D19" "46.590.99"
Set-TextValue $ws "E19" "  -0.46%  "
Set-TextValue $ws "E20" "  +1.41%  "
Set-TextValue $ws "D21" "13.06"
Set-TextValue $ws "E21" "  -5.23%  "
Set-TextValue $ws "D22" "6.79"
Set-TextValue $ws "E22" "  +3.13%  "
Set-TextValue $ws "E23" "  +2.67%  "
Set-TextValue $ws "D24" "276.20"
Set-TextValue $ws "E24" "  +8.39%  "
Set-TextValue $ws "D25" "3.04"
Set-TextValue $ws "E25" "  +1.95%  "
Set-TextValue $ws "D26" "2.22"
Set-TextValue $ws "E26" "  +4.98%  "
Set-TextValue $ws "D27" "29.55"
Set-TextValue $ws "E27" "  +13.63%  "
Set-TextValue $ws "D28" "1.00"
Set-TextValue $ws "E28" "  +0.19%  "
Set-TextValue $ws "D29" "4.03"
Set-TextValue $ws "E29" "  -1.53%  "
Set-TextValue $ws "E30" "  +2.30%  "
Set-TextValue $ws "D31" "38.67"
Set-TextValue $ws "E31" "  -5.78%  "
Set-TextValue $ws "D32" "2.23"
Set-TextValue $ws "E32" "  -2.67%  "
Set-TextValue $ws "D33" "6.45"
Set-TextValue $ws "E33" "  +8.44%  "
Set-TextValue $ws "D34" "3.65"
Set-TextValue $ws "E34" "  -3.51%  "
Set-TextValue $ws "D35" "2.26"
Set-TextValue $ws "E35" "  +2.00%  "
Set-TextValue $ws "E36" "  -3.27%  "
Set-TextValue $ws "E37" "  -0.37%  "
Set-TextValue $ws "D38" "152.40"
Set-TextValue $ws "E38" "  +2.37%  "
Set-TextValue $ws "E39" "  +0.45%  "
Set-TextValue $ws "E40" "  +2.20%  "
Set-TextValue $ws "D41" "24.14"
Set-TextValue $ws "E41" "  +36.78%  "
Set-TextValue $ws "D42" "15.97"
Set-TextValue $ws "E42" "  -2.30%  "
Set-TextValue $ws "E43" "  +1.08%  "
Set-TextValue $ws "D44" "3.62"
Set-TextValue $ws "E44" "  +1.51%  "
Set-TextValue $ws "D45" "4.07"
Set-TextValue $ws "E45" "  -4.01%  "
Set-TextValue $ws "D46" "2.141.55"
Set-TextValue $ws "E46" "  +6.03%  "
Set-TextValue $ws "E47" "  -0.14%  "
Set-TextValue $ws "D48" "95.03"
Set-TextValue $ws "E48" "  +2.10%  "
Set-TextValue $ws "D49" "9.52"
Set-TextValue $ws "E49" "  +8.10%  "
Set-TextValue $ws "D50" "109.70"
Set-TextValue $ws "E50" "  +2.12%  "
Set-TextValue $ws "E51" "  -4.14%  "
